$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 7 — appended observation record (numeric columns)
$ws.Range("A7").Value = 112182932
$ws.Range("B7").Value = 77515
$ws.Range("E7").Value = 6425
$ws.Range("Q7").Value = 527228.9315337478
$ws.Range("R7").Value = 6908168.611324663
$ws.Range("S7").Value = 5

# Text columns (stored as inline/shared strings in the target)
$ws.Range("C7").Value = "Ovaliderad"
$ws.Range("D7").Value = "NT"
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "1"
$ws.Range("I7").Style = "Normal"
$ws.Range("P7").Value = "Nylandet, Hls"
$ws.Range("T7").Value = "Gävleborg"
$ws.Range("U7").Value = "Ljusdal"
$ws.Range("V7").Value = "Hälsingland"
$ws.Range("W7").Value = "Ramsjö"
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value = "2023-07-06"
$ws.Range("Y7").Style = "Normal"
$ws.Range("Z7").Value = "00:00"
$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = "2023-07-06"
$ws.Range("AA7").Style = "Normal"
$ws.Range("AB7").Value = "00:00"
$ws.Range("AC7").Value = "Påträffad under Sveaskogs naturvärdesinventering"
$ws.Range("AW7").Value = "Mimmi Persson"
$ws.Range("AX7").Value = "Mimmi Persson"

# Boolean columns
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false

# Empty-but-present text cells (leading apostrophe collapses to an empty
# string cell, matching the other rows' blank inlineStr cells)
$ws.Range("AT7").Value = "'"
$ws.Range("AT7").Style = "Normal"
$ws.Range("AY7").Value = "'"
$ws.Range("AY7").Style = "Normal"
